$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 168501.83
$ws.Range("I4").Value = 168501.83
$ws.Range("K4").Value = 168501.83
$ws.Range("M4").Value = -168387.83
$ws.Range("H18").Value = 6496.5
$ws.Range("J18").Value = 6494.5
$ws.Range("L18").Value = 6494.5
$ws.Range("N18").Value = -7062.5
$ws.Range("H103").Value = 1444.6364
$ws.Range("I103").Value = 332.5
$ws.Range("K103").Value = 997.5
$ws.Range("M103").Value = -411.5
$ws.Range("H133").Value = 170494.86
$ws.Range("J133").Value = 182792.5
$ws.Range("L133").Value = 182792.5
$ws.Range("N133").Value = -192912.5
$ws.Range("H138").Value = 2144.2058
$ws.Range("I138").Value = 1681.4375
$ws.Range("K138").Value = 5044.3125
$ws.Range("M138").Value = 95.6875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3230.3975
$ws.Range("I32").Value = 2663.145
$ws.Range("K32").Value = 2663.145
$ws.Range("M32").Value = -2376.145
$ws.Range("H74").Value = 18305.66
$ws.Range("I74").Value = 1127.1177
$ws.Range("J74").Value = 49046.21
$ws.Range("K74").Value = 1127.1177
$ws.Range("L74").Value = 49046.21
$ws.Range("M74").Value = -253.1177
$ws.Range("N74").Value = -50794.21
$ws.Range("H77").Value = 18305.66
$ws.Range("I77").Value = 1127.1177
$ws.Range("J77").Value = 49046.21
$ws.Range("K77").Value = 5635.5885
$ws.Range("L77").Value = 245231.05
$ws.Range("M77").Value = -1267.5885
$ws.Range("N77").Value = -253967.05
$ws.Range("H122").Value = 50076.473
$ws.Range("I122").Value = 2477.1667
$ws.Range("K122").Value = 7431.500100000001
$ws.Range("M122").Value = -4981.500100000001
$ws.Range("H132").Value = 2631.5
$ws.Range("I132").Value = 2553.83
$ws.Range("K132").Value = 7661.49
$ws.Range("M132").Value = -5131.49

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2392.5386
$ws.Range("I99").Value = 2642.2222
$ws.Range("J99").Value = 1830.75
$ws.Range("K99").Value = 2642.2222
$ws.Range("L99").Value = 1830.75
$ws.Range("M99").Value = -1144.2222
$ws.Range("N99").Value = -4826.75
$ws.Range("H105").Value = 2399.3333
$ws.Range("I105").Value = 2643.6667
$ws.Range("J105").Value = 1910.6666
$ws.Range("K105").Value = 2643.6667
$ws.Range("L105").Value = 1910.6666
$ws.Range("M105").Value = -896.6667000000002
$ws.Range("N105").Value = -5404.6666
$ws.Range("H107").Value = 6477
$ws.Range("I107").Value = 6477
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 6477
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4557
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2012.9081
$ws.Range("I31").Value = 1338.5526
$ws.Range("J31").Value = 2535.8774
$ws.Range("K31").Value = 1338.5526
$ws.Range("L31").Value = 2535.8774
$ws.Range("M31").Value = -1043.5526
$ws.Range("N31").Value = -3125.8774
$ws.Range("H34").Value = 2012.9081
$ws.Range("I34").Value = 1338.5526
$ws.Range("J34").Value = 2535.8774
$ws.Range("K34").Value = 1338.5526
$ws.Range("L34").Value = 2535.8774
$ws.Range("M34").Value = -1136.5526
$ws.Range("N34").Value = -2939.8774
$ws.Range("H86").Value = 11101.182
$ws.Range("I86").Value = 8822.308000000001
$ws.Range("K86").Value = 8822.308000000001
$ws.Range("M86").Value = -7699.308000000001
$ws.Range("H89").Value = 11101.182
$ws.Range("I89").Value = 8822.308000000001
$ws.Range("K89").Value = 44111.54000000001
$ws.Range("M89").Value = -38495.54000000001
$ws.Range("H105").Value = 1520.8
$ws.Range("I105").Value = 1368.6666
$ws.Range("J105").Value = 1749
$ws.Range("K105").Value = 1368.6666
$ws.Range("L105").Value = 1749
$ws.Range("M105").Value = 378.3334
$ws.Range("N105").Value = -5243
$ws.Range("H107").Value = 100092760
$ws.Range("I107").Value = 166817890
$ws.Range("K107").Value = 166817890
$ws.Range("M107").Value = -166815970
$ws.Range("H122").Value = 3254.1667
$ws.Range("I122").Value = 2831.25
$ws.Range("J122").Value = 4100
$ws.Range("K122").Value = 8493.75
$ws.Range("L122").Value = 12300
$ws.Range("M122").Value = -6043.75
$ws.Range("N122").Value = -17200

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1545.5385
$ws.Range("I5").Value = 582
$ws.Range("K5").Value = 1746
$ws.Range("M5").Value = -1634
$ws.Range("H135").Value = 1545.5385
$ws.Range("I135").Value = 582
$ws.Range("K135").Value = 5238
$ws.Range("M135").Value = -2703

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4247.722
$ws.Range("I40").Value = 4411.6875
$ws.Range("K40").Value = 4411.6875
$ws.Range("M40").Value = -4275.6875
$ws.Range("H68").Value = 4499.5454
$ws.Range("I68").Value = 2277.2222
$ws.Range("J68").Value = 14500
$ws.Range("K68").Value = 2277.2222
$ws.Range("L68").Value = 14500
$ws.Range("M68").Value = -1528.2222
$ws.Range("N68").Value = -15998
$ws.Range("H71").Value = 4499.5454
$ws.Range("I71").Value = 2277.2222
$ws.Range("J71").Value = 14500
$ws.Range("K71").Value = 11386.111
$ws.Range("L71").Value = 72500
$ws.Range("M71").Value = -7642.111000000001
$ws.Range("N71").Value = -79988
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13886.444
$ws.Range("I62").Value = 11998
$ws.Range("J62").Value = 14426
$ws.Range("K62").Value = 11998
$ws.Range("L62").Value = 14426
$ws.Range("M62").Value = -11374
$ws.Range("N62").Value = -15674
$ws.Range("H65").Value = 13886.444
$ws.Range("I65").Value = 11998
$ws.Range("J65").Value = 14426
$ws.Range("K65").Value = 59990
$ws.Range("L65").Value = 72130
$ws.Range("M65").Value = -56870
$ws.Range("N65").Value = -78370
$ws.Range("H75").Value = 126249.5
$ws.Range("I75").Value = 130000
$ws.Range("J75").Value = 124999.336
$ws.Range("K75").Value = 130000
$ws.Range("L75").Value = 124999.336
$ws.Range("M75").Value = -129064
$ws.Range("N75").Value = -126871.336
$ws.Range("H78").Value = 126249.5
$ws.Range("I78").Value = 130000
$ws.Range("J78").Value = 124999.336
$ws.Range("K78").Value = 390000
$ws.Range("L78").Value = 374998.008
$ws.Range("M78").Value = -385320
$ws.Range("N78").Value = -384358.008
$ws.Range("H125").Value = 93333
$ws.Range("J125").Value = 93333
$ws.Range("L125").Value = 93333
$ws.Range("N125").Value = -103173
$ws.Range("H132").Value = 3055.9546
$ws.Range("I132").Value = 3225.842
$ws.Range("J132").Value = 1980
$ws.Range("K132").Value = 9677.526
$ws.Range("L132").Value = 5940
$ws.Range("M132").Value = -7147.526
$ws.Range("N132").Value = -11000
